$wb = $excel.ActiveWorkbook

$oldId = "5d5eae22-f549-4ce7-bf34-414228d83089"
$newId = "7e1fce7d-ad49-4108-ad0f-5ebb7df888b1"
$newHash = "4ee12fb2bfee37fb67194804ee5a07cd4761648a"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-09-05 15:14:36"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c25edf8fe20d7e6c7d8c7d087c7f176687fae871/e2e/$newId.md", "", "", "e2e\$newId.md") | Out-Null

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-05 15:14:32"
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c25edf8fe20d7e6c7d8c7d087c7f176687fae871/e2e/$newId.md", "", "", "$newId.md") | Out-Null

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-05 15:14:36"
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c25edf8fe20d7e6c7d8c7d087c7f176687fae871/e2e/$newId.md", "", "", "$newId.md") | Out-Null

# --- Column widths for I/J on zh-cn and de-de ---
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
